$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = -0.0703
$ws.Range("K6").Value = -0.0731
$ws.Range("L6").Value = -0.055
$ws.Range("M6").Value = -0.0487
$ws.Range("N6").Value = -0.0394
$ws.Range("O6").Value = -0.0275
$ws.Range("P6").Value = -0.0246

$ws.Range("N8").Value = 0.1478
$ws.Range("O8").Value = 0.1643
$ws.Range("P8").Value = 0.2998

$ws.Range("I15").Value = -2.0918
$ws.Range("J15").Value = -2.087
$ws.Range("K15").Value = -2.9915
$ws.Range("L15").Value = -1.2599
$ws.Range("M15").Value = -0.9552
$ws.Range("N15").Value = -0.4207
$ws.Range("O15").Value = -0.5797
$ws.Range("P15").Value = -0.4185

$ws.Range("I22").Value = -0.0172
$ws.Range("J22").Value = -0.0195
$ws.Range("K22").Value = -0.0206
$ws.Range("L22").Value = -0.0122
$ws.Range("M22").Value = -0.0104
$ws.Range("N22").Value = -0.0067
$ws.Range("O22").Value = -0.0037
$ws.Range("P22").Value = -0.002

$ws.Range("I24").Value = -0.1418
$ws.Range("J24").Value = -0.0963
$ws.Range("K24").Value = -0.0672
$ws.Range("L24").Value = -0.062
$ws.Range("M24").Value = -0.0451
$ws.Range("N24").Value = -0.0348
$ws.Range("O24").Value = -0.0272
$ws.Range("P24").Value = -0.0009

$ws.Range("J33").Value = -0.0268
$ws.Range("K33").Value = -0.0265
$ws.Range("L33").Value = -0.0262
$ws.Range("M33").Value = -0.0258
$ws.Range("N33").Value = -0.0254
$ws.Range("O33").Value = -0.025
$ws.Range("P33").Value = -0.0247

$ws.Range("N35").Value = 0.0037
$ws.Range("O35").Value = 0.011
$ws.Range("P35").Value = 0.0054

$ws.Range("I42").Value = 0.2094
$ws.Range("J42").Value = 0.1077
$ws.Range("K42").Value = 0.0725
$ws.Range("L42").Value = 0.0262
$ws.Range("M42").Value = -0.073
$ws.Range("N42").Value = -0.0583
$ws.Range("O42").Value = -0.0495
$ws.Range("P42").Value = -0.0496

$ws.Range("I49").Value = -0.0052
$ws.Range("J49").Value = -0.005
$ws.Range("K49").Value = -0.0049
$ws.Range("L49").Value = -0.0047
$ws.Range("M49").Value = -0.0046
$ws.Range("N49").Value = -0.0045
$ws.Range("O49").Value = -0.0044
$ws.Range("P49").Value = -0.0043

$ws.Range("I51").Value = -0.01
$ws.Range("J51").Value = -0.005
$ws.Range("K51").Value = -0.0049
$ws.Range("L51").Value = -0.0048
$ws.Range("M51").Value = -0.0048
$ws.Range("N51").Value = -0.0047
$ws.Range("O51").Value = -0.0046
$ws.Range("P51").Value = 0
